# Apply crypto price/volume updates to Sheet1, preserving each cell's
# original text data type (the sheet stores Price/Volume as text, not
# numbers, so number-looking strings must not be auto-coerced by Excel).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '60.132.24'
$ws.Range("E2").Value = '  +3.79%  '
$ws.Range("D3").Value = '3.201.66'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '538.47'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.70%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '145.46'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +4.74%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.999'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("E8").Value = '  +1.79%  '
$ws.Range("E9").Value = '  -0.19%  '
$ws.Range("E10").Value = '  +4.93%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.429'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +2.56%  '
$ws.Range("D12").Value = '3.749.96'
$ws.Range("E12").Value = '  +2.40%  '
$ws.Range("E13").Value = '  -0.92%  '
$ws.Range("B14").Value = 'Avalanche'
$ws.Range("C14").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '26.19'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.83%  '
$ws.Range("B15").Value = 'ShibaInu'
$ws.Range("C15").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0000175'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +3.82%  '
$ws.Range("D16").Value = '60.130.89'
$ws.Range("E16").Value = '  +3.68%  '
$ws.Range("D17").Value = '3.217.33'
$ws.Range("E17").Value = '  +2.91%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '6.22'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.06%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.12'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.04%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '8.38'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +2.86%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '382.24'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.90%  '
$ws.Range("E22").Value = '  +0.18%  '
$ws.Range("E23").Value = '  +3.63%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '70.11'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.39%  '
$ws.Range("E25").Value = '  +3.22%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '8.83'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +14.00%  '
$ws.Range("E27").Value = '  -0.07%  '
$ws.Range("D28").Value = '0.0₃0906'
$ws.Range("E28").Value = '  +1.97%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.92'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +2.13%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '22.50'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +3.66%  '
$ws.Range("E31").Value = '  +5.37%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.18'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.06%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.22'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +3.73%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.51'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +4.92%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '156.76'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -2.52%  '
$ws.Range("E36").Value = '  +2.37%  '
$ws.Range("D37").Value = '2.776.32'
$ws.Range("E37").Value = '  +8.96%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '25.71'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.87%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0711'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +5.19%  '
$ws.Range("E40").Value = '  +1.09%  '
$ws.Range("E41").Value = '  +1.65%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.731'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +4.52%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '39.73'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +2.76%  '
$ws.Range("E44").Value = '  +5.88%  '
$ws.Range("D45").Value = '3.244.31'
$ws.Range("E45").Value = '  +2.52%  '
$ws.Range("E46").Value = '  +3.13%  '
$ws.Range("E47").Value = '  -0.15%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.100'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.66%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '20.67'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +3.17%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.795'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +5.94%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.00'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.05%  '
